$p = $ppt.ActivePresentation

# The deck originally bundled slides for several careers (Electrica,
# Mecatronica, Sistemas Computacionales). The generator now emits one
# .pptx per career, so for "INGENIERIA ELECTRICA" we keep only the
# first 3 slides and drop the trailing slides (and their notes pages)
# that belong to other careers.
while ($p.Slides.Count -gt 3) {
    $p.Slides.Item($p.Slides.Count).Delete()
}
